$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3288.75
$ws.Range("J112").Value = 3288.75
$ws.Range("L112").Value = 9866.25
$ws.Range("N112").Value = -12082.25
$ws.Range("H116").Value = 2904.7222
$ws.Range("I116").Value = 2814.9167
$ws.Range("K116").Value = 2814.9167
$ws.Range("M116").Value = 627.0832999999998
$ws.Range("H129").Value = 932
$ws.Range("I129").Value = 795.7143
$ws.Range("J129").Value = 988.1177
$ws.Range("K129").Value = 2387.1429
$ws.Range("L129").Value = 2964.3531
$ws.Range("M129").Value = 2612.8571
$ws.Range("N129").Value = -12964.3531
$ws.Range("H132").Value = 5958003
$ws.Range("I132").Value = 6947960.5
$ws.Range("J132").Value = 18260.625
$ws.Range("K132").Value = 20843881.5
$ws.Range("L132").Value = 54781.875
$ws.Range("M132").Value = -20841351.5
$ws.Range("N132").Value = -59841.875
$ws.Range("H137").Value = 1566.4138
$ws.Range("I137").Value = 1429.5555
$ws.Range("J137").Value = 1790.3636
$ws.Range("K137").Value = 4288.666499999999
$ws.Range("L137").Value = 5371.0908
$ws.Range("M137").Value = -1738.666499999999
$ws.Range("N137").Value = -10471.0908
$ws.Range("H138").Value = 2989.3552
$ws.Range("J138").Value = 2993.2222
$ws.Range("L138").Value = 8979.6666
$ws.Range("N138").Value = -19259.6666

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12727.714
$ws.Range("I32").Value = 10718.62
$ws.Range("J32").Value = 14820.5205
$ws.Range("K32").Value = 10718.62
$ws.Range("L32").Value = 14820.5205
$ws.Range("M32").Value = -10431.62
$ws.Range("N32").Value = -15394.5205
$ws.Range("H61").Value = 111112340
$ws.Range("I61").Value = 166667500
$ws.Range("J61").Value = 2004.6666
$ws.Range("K61").Value = 166667500
$ws.Range("L61").Value = 2004.6666
$ws.Range("M61").Value = -166667288
$ws.Range("N61").Value = -2428.6666
$ws.Range("H63").Value = 2590.0667
$ws.Range("I63").Value = 2463.8333
$ws.Range("J63").Value = 3095
$ws.Range("K63").Value = 2463.8333
$ws.Range("L63").Value = 3095
$ws.Range("M63").Value = -1777.8333
$ws.Range("N63").Value = -4467
$ws.Range("H66").Value = 2590.0667
$ws.Range("I66").Value = 2463.8333
$ws.Range("J66").Value = 3095
$ws.Range("K66").Value = 12319.1665
$ws.Range("L66").Value = 15475
$ws.Range("M66").Value = -8887.166499999999
$ws.Range("N66").Value = -22339
$ws.Range("H136").Value = 111112340
$ws.Range("I136").Value = 166667500
$ws.Range("J136").Value = 2004.6666
$ws.Range("K136").Value = 500002500
$ws.Range("L136").Value = 6013.9998
$ws.Range("M136").Value = -499999950
$ws.Range("N136").Value = -11113.9998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7067.0625
$ws.Range("I134").Value = 876.6429000000001
$ws.Range("J134").Value = 50400
$ws.Range("K134").Value = 2629.9287
$ws.Range("L134").Value = 151200
$ws.Range("M134").Value = -94.92870000000039
$ws.Range("N134").Value = -156270

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 211.91667
$ws.Range("I7").Value = 160.33333
$ws.Range("K7").Value = 160.33333
$ws.Range("M7").Value = -47.33332999999999
$ws.Range("H16").Value = 76924340
$ws.Range("I16").Value = 76924340
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 76924340
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -76924053
$ws.Range("H58").Value = 4766.759
$ws.Range("I58").Value = 752.4375
$ws.Range("J58").Value = 9707.462
$ws.Range("K58").Value = 752.4375
$ws.Range("L58").Value = 9707.462
$ws.Range("M58").Value = -549.4375
$ws.Range("N58").Value = -10113.462
$ws.Range("H113").Value = 76924340
$ws.Range("I113").Value = 76924340
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 76924340
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -76922170
$ws.Range("H114").Value = 31732.8
$ws.Range("J114").Value = 31732.8
$ws.Range("L114").Value = 31732.8
$ws.Range("N114").Value = -40410.8
$ws.Range("H136").Value = 4766.759
$ws.Range("I136").Value = 752.4375
$ws.Range("J136").Value = 9707.462
$ws.Range("K136").Value = 2257.3125
$ws.Range("L136").Value = 29122.386
$ws.Range("M136").Value = 292.6875
$ws.Range("N136").Value = -34222.386
$ws.Range("H141").Value = 264884.2
$ws.Range("J141").Value = 264884.2
$ws.Range("L141").Value = 264884.2
$ws.Range("N141").Value = -275244.2

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 975
$ws.Range("J54").Value = 975
$ws.Range("L54").Value = 2925
$ws.Range("N54").Value = -4043
$ws.Range("H107").Value = 6468.773
$ws.Range("I107").Value = 490
$ws.Range("K107").Value = 1470
$ws.Range("M107").Value = 450
$ws.Range("H122").Value = 932.6667
$ws.Range("J122").Value = 1010
$ws.Range("L122").Value = 9090
$ws.Range("N122").Value = -13990
$ws.Range("H131").Value = 18897178
$ws.Range("I131").Value = 200000380
$ws.Range("J131").Value = 32262.271
$ws.Range("K131").Value = 600001140
$ws.Range("L131").Value = 96786.81299999999
$ws.Range("M131").Value = -599996100
$ws.Range("N131").Value = -106866.813
$ws.Range("H137").Value = 31257506
$ws.Range("J137").Value = 10876.071
$ws.Range("L137").Value = 32628.213
$ws.Range("N137").Value = -42828.213

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 9615484
$ws.Range("I107").Value = 9615484
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 9615484
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -9613564
$ws.Range("H132").Value = 9383.166999999999
$ws.Range("I132").Value = 11344.923
$ws.Range("J132").Value = 4282.6
$ws.Range("K132").Value = 34034.769
$ws.Range("L132").Value = 12847.8
$ws.Range("M132").Value = -31504.769
$ws.Range("N132").Value = -17907.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = 0
$ws.Range("H22").Value = 1800.7142
$ws.Range("J22").Value = 1340.8
$ws.Range("L22").Value = 1340.8
$ws.Range("N22").Value = -1930.8
$ws.Range("H27").Value = 1800.7142
$ws.Range("J27").Value = 1340.8
$ws.Range("L27").Value = 1340.8
$ws.Range("N27").Value = -1554.8
$ws.Range("H100").Value = 1542.4286
$ws.Range("I100").Value = 1199.25
$ws.Range("K100").Value = 1199.25
$ws.Range("M100").Value = -658.25
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = 0

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 716
$ws.Range("I113").Value = 324
$ws.Range("K113").Value = 972
$ws.Range("M113").Value = 1198
